$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the stray "asd" test entry that lived on row 35 (E35/F35).
# Clearing both cells empties the row entirely, and the now-unreferenced
# "asd" shared string drops out of the shared-strings table automatically.
$ws.Range("E35").Value = $null
$ws.Range("F35").Value = $null

# Row 11 ("Problem Set 2" session) actually ran later than first logged -
# update its end time (columns D/E were not filled in before).
$ws.Range("B11").Value = 20
$ws.Range("C11").Value = 0
$ws.Range("D11").Value = 21
$ws.Range("E11").Value = 40

# Log the new "lecture 12 finger exercise" session as row 12.
$ws.Range("A12").Value = 45812
$ws.Range("A11").Copy()
$ws.Range("A12").PasteSpecial(-4122)
$ws.Range("B12").Value = 21
$ws.Range("C12").Value = 10
$ws.Range("F12").Value = "CS introduction Lecture 12"

# Leave the selection where the author left off editing.
$ws.Range("C12").Select()
